$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2272727272727273
$ws.Range("C2").Value = 0.4734848484848485
$ws.Range("J2").Value = 0.01893939393939394
$ws.Range("P2").Value = 0.1515151515151515
$ws.Range("S2").Value = 0.1287878787878788
$ws.Range("C3").Value = 0.04615384615384616
$ws.Range("J3").Value = 0.03076923076923077
$ws.Range("P3").Value = 0.7153846153846154
$ws.Range("S3").Value = 0.2076923076923077
$ws.Range("J4").Value = 0.04545454545454546
$ws.Range("P4").Value = 0.6818181818181818
$ws.Range("S4").Value = 0.2727272727272727
$ws.Range("B6").Value = 0.06770833333333333
$ws.Range("D6").Value = 0.01041666666666667
$ws.Range("F6").Value = 0.0625
$ws.Range("J6").Value = 0.1822916666666667
$ws.Range("O6").Value = 0.015625
$ws.Range("Q6").Value = 0.1822916666666667
$ws.Range("R6").Value = 0.1041666666666667
$ws.Range("S6").Value = 0.375
$ws.Range("B7").Value = 0.1141304347826087
$ws.Range("D7").Value = 0.03804347826086957
$ws.Range("F7").Value = 0.05434782608695652
$ws.Range("J7").Value = 0.1684782608695652
$ws.Range("O7").Value = 0.01630434782608696
$ws.Range("Q7").Value = 0.1630434782608696
$ws.Range("R7").Value = 0.09239130434782608
$ws.Range("S7").Value = 0.3532608695652174
$ws.Range("B8").Value = 0.07032967032967033
$ws.Range("D8").Value = 0.01978021978021978
$ws.Range("E8").Value = 0.002197802197802198
$ws.Range("F8").Value = 0.07252747252747253
$ws.Range("J8").Value = 0.1384615384615385
$ws.Range("O8").Value = 0.01318681318681319
$ws.Range("Q8").Value = 0.1846153846153846
$ws.Range("R8").Value = 0.1186813186813187
$ws.Range("S8").Value = 0.3802197802197802
$ws.Range("B9").Value = 0.08641975308641975
$ws.Range("D9").Value = 0.01851851851851852
$ws.Range("F9").Value = 0.09259259259259259
$ws.Range("J9").Value = 0.1296296296296296
$ws.Range("O9").Value = 0.01234567901234568
$ws.Range("Q9").Value = 0.154320987654321
$ws.Range("R9").Value = 0.1172839506172839
$ws.Range("S9").Value = 0.3888888888888889
$ws.Range("B10").Value = 0.09790794979079498
$ws.Range("D10").Value = 0.02092050209205021
$ws.Range("E10").Value = 0.001673640167364017
$ws.Range("F10").Value = 0.06694560669456066
$ws.Range("J10").Value = 0.1171548117154812
$ws.Range("O10").Value = 0.02175732217573222
$ws.Range("Q10").Value = 0.2150627615062762
$ws.Range("R10").Value = 0.09205020920502092
$ws.Range("S10").Value = 0.3665271966527197
$ws.Range("G11").Value = 0.1081081081081081
$ws.Range("J11").Value = 0.1013513513513514
$ws.Range("K11").Value = 0.1993243243243243
$ws.Range("L11").Value = 0.5641891891891891
$ws.Range("S11").Value = 0.02702702702702703
$ws.Range("G12").Value = 0.7485380116959064
$ws.Range("J12").Value = 0.1695906432748538
$ws.Range("K12").Value = 0.01169590643274854
$ws.Range("L12").Value = 0.02923976608187134
$ws.Range("S12").Value = 0.04093567251461988
$ws.Range("G13").Value = 0.6595744680851063
$ws.Range("J13").Value = 0.2978723404255319
$ws.Range("S13").Value = 0.0425531914893617
$ws.Range("F15").Value = 0.00966183574879227
$ws.Range("H15").Value = 0.1884057971014493
$ws.Range("I15").Value = 0.08695652173913043
$ws.Range("J15").Value = 0.3671497584541063
$ws.Range("K15").Value = 0.03864734299516908
$ws.Range("M15").Value = 0.00966183574879227
$ws.Range("O15").Value = 0.03864734299516908
$ws.Range("S15").Value = 0.2608695652173913
$ws.Range("F16").Value = 0.01290322580645161
$ws.Range("H16").Value = 0.1741935483870968
$ws.Range("I16").Value = 0.09032258064516129
$ws.Range("J16").Value = 0.3935483870967742
$ws.Range("K16").Value = 0.09032258064516129
$ws.Range("M16").Value = 0.02580645161290323
$ws.Range("O16").Value = 0.07096774193548387
$ws.Range("S16").Value = 0.1419354838709677
$ws.Range("F17").Value = 0.007009345794392523
$ws.Range("H17").Value = 0.1799065420560748
$ws.Range("I17").Value = 0.07943925233644859
$ws.Range("J17").Value = 0.4345794392523364
$ws.Range("K17").Value = 0.08644859813084112
$ws.Range("M17").Value = 0.02803738317757009
$ws.Range("O17").Value = 0.0630841121495327
$ws.Range("S17").Value = 0.1214953271028037
$ws.Range("F18").Value = 0.004524886877828055
$ws.Range("H18").Value = 0.2081447963800905
$ws.Range("I18").Value = 0.07239819004524888
$ws.Range("J18").Value = 0.4660633484162896
$ws.Range("K18").Value = 0.09049773755656108
$ws.Range("M18").Value = 0.009049773755656109
$ws.Range("N18").Value = 0.004524886877828055
$ws.Range("O18").Value = 0.04524886877828054
$ws.Range("S18").Value = 0.09954751131221719
$ws.Range("F19").Value = 0.01016949152542373
$ws.Range("H19").Value = 0.2245762711864407
$ws.Range("I19").Value = 0.06864406779661017
$ws.Range("J19").Value = 0.35
$ws.Range("K19").Value = 0.1338983050847458
$ws.Range("M19").Value = 0.02372881355932203
$ws.Range("N19").Value = 0.000847457627118644
$ws.Range("O19").Value = 0.06694915254237288
$ws.Range("S19").Value = 0.1211864406779661
